# Update odds values on Sheet1 to match the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67

# Row 7 updates
$ws.Range("G7").Value = 4.05
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 2.37
$ws.Range("N7").Value = 6.9
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.05
$ws.Range("T7").Value = 2.75
$ws.Range("W7").Value = 11
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 14
$ws.Range("Z7").Value = 70
$ws.Range("AB7").Value = 50
$ws.Range("AC7").Value = 6.9
$ws.Range("AF7").Value = 80
$ws.Range("AL7").Value = 28
$ws.Range("AN7").Value = 6
$ws.Range("AO7").Value = 24
$ws.Range("AQ7").Value = 150
$ws.Range("AT7").Value = 2.75
$ws.Range("AX7").Value = 8.75
